$wb = $excel.ActiveWorkbook

# --- Update "tasas" sheet numeric rates ---
$tasas = $wb.Worksheets.Item("tasas")

$tasas.Range("N10").Value = 132.29
$tasas.Range("O10").Value = 4040

$tasas.Range("N12").Value = 4040
$tasas.Range("O12").Value = 129

# --- Update "Hoja1" conversion summary text ---
$hoja1 = $wb.Worksheets.Item("Hoja1")

$oldLine1 = [string][char]0x2705 + " 1000 Bs = 7.52 = 30338.13 pesos"
$newLine1 = [string][char]0x2705 + " 1000 Bs = 7.56 = 30538.97 pesos"

$oldLine2 = [string][char]0x2705 + " 30338.13 pesos = 7.49 = 943.85 Bs"
$newLine2 = [string][char]0x2705 + " 30538.97 pesos = 7.56 = 975.13 Bs"

$current = $hoja1.Range("A1").Value()
$updated = $current.Replace($oldLine1, $newLine1)
$updated = $updated.Replace($oldLine2, $newLine2)
$hoja1.Range("A1").Value = $updated
